$d = $word.ActiveDocument

# Locate the run that currently reads
#   "Nom du programme Python :  READ_CSV_FINAL.py "
# and split it into five runs (same bold formatting) whose texts are:
#   "Nom du programme Python :  READ_CSV_FINAL", "_", " ",
#   "FRIDHI_Ilies_LEGER_Maureen", ".py "

$old = "Nom du programme Python :  READ_CSV_FINAL.py "

$target = $d.Content
$target.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $target.Find.Found) {
    throw "Could not find the target run text to split."
}

# Remember where the (soon to be removed) paragraph mark of the
# now-empty host paragraph will end up, so the two paragraphs can be
# re-merged after the XML insertion below.
$paraIndexBefore = $d.Paragraphs.Count

# Clear the existing run text; InsertXML below will repopulate the
# paragraph with the newly split runs.
$target.Delete()

$rPr = '<w:rPr><w:b/><w:bCs/></w:rPr>'
$nsW = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$runs =
    ('<w:r ' + $nsW + '>' + $rPr + '<w:t>Nom du programme Python :  READ_CSV_FINAL</w:t></w:r>') +
    ('<w:r ' + $nsW + '>' + $rPr + '<w:t>_</w:t></w:r>') +
    ('<w:r ' + $nsW + '>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>') +
    ('<w:r ' + $nsW + '>' + $rPr + '<w:t>FRIDHI_Ilies_LEGER_Maureen</w:t></w:r>') +
    ('<w:r ' + $nsW + '>' + $rPr + '<w:t>.py </w:t></w:r>')

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $runs + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($payload)

# InsertXML inserted the five runs as a brand-new paragraph immediately
# before the original (now empty) host paragraph, which still carries the
# original paragraph-level formatting/rsid attributes. Re-merge the two
# paragraphs by deleting the paragraph mark between them so the original
# paragraph's properties apply to the new runs, matching a simple in-place
# run split rather than leaving a stray empty paragraph behind.
$newHostPara = $d.Paragraphs.Item($paraIndexBefore)
$markStart = $newHostPara.Range.End - 1
$d.Range($markStart, $markStart + 1).Delete()
